$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add "history" (G1) and "balance" (H1) columns ---

# G1 "history" - reuse the same bold/centered/bordered header look as the other headers
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value = "history"

# H1 "balance" - same bold/centered look, but only left+right borders (no top/bottom)
$ws.Range("F1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "balance"
$ws.Range("H1").Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> xlLineStyleNone
$ws.Range("H1").Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> xlLineStyleNone

# --- Row 2: fix the "history" cell content, add the new "balance" value ---

$ws.Range("G2").Style = "Normal"
$ws.Range("G2").Value = "Missing|injasdf"
$ws.Range("H2").Value = 0

# --- Row 3: brand new member record ---

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "735554"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "morteza"

$ws.Range("C3").Value = "pashaei"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "926010932"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0926010932"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "03/02/2023"
$ws.Range("F3").Style = "Normal"

$ws.Range("G3").Value = "Missing"

$ws.Range("H3").Value = 0

# --- Selection matches the saved file ---
$ws.Range("G1").Select()
